{"js": "// Office.js (Word JavaScript API) edit script.\n// Body: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) Date line\nitems[0].insertText(\"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 04.07.24:\u26a1\ufe0f\ud83d\ude80\", \"Replace\");\n\n// 2) Title\nitems[1].insertText(\n  \"How Do Large Language Models Acquire Factual Knowledge During Pretraining?\",\n  \"Replace\"\n);\n\n// 3) Body paragraph 1\nitems[2].insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 \u05e0\u05d5\u05e9\u05d0 \u05de\u05ea\u05d9 \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d0\u05e9\u05db\u05e8\u05d4 \u05e8\u05d5\u05db\u05e9\u05d9\u05dd \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 (\u05dc\u05de\u05e9\u05dc \u05e9\u05e2\u05d9\u05e8 \u05d1\u05d9\u05e8\u05d4 \u05e9\u05dc \u05e6\u05e8\u05e4\u05ea \u05d4\u05d9\u05d0 \u05e4\u05e8\u05d9\u05e1) \u05d1\u05de\u05d4\u05dc\u05da \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05e7\u05d3\u05d9\u05dd. \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05d2\u05dd \u05d1\u05d5\u05d3\u05e7 \u05db\u05de\u05d4 \u05d6\u05de\u05df \u05dc\u05d5\u05e7\u05d7 \u05dc\u05e9\u05db\u05d5\u05d7 \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9. \u05d0\u05d5\u05e7\u05d9\u05d9, \u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05e9\u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05e9\u05dc\u05e0\u05d5 \u05e2\u05dd \u05d0\u05d7\u05ea \u05d4\u05e6\u05d5\u05e8\u05d5\u05ea \u05e9\u05dc \u05de\u05e9\u05e4\u05d7\u05ea \u05de\u05d5\u05e8\u05d3 \u05d4\u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 (gradient descent \u05d0\u05d5 GD). \u05d1\u05d3\u05f4\u05db \u05d3\u05d5\u05d2\u05de\u05d9\u05dd \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d4\u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc\u05e0\u05d5 (\u05de\u05d9\u05e0\u05d9-\u05d1\u05d0\u05e5') \u05d5\u05de\u05d6\u05d9\u05d6\u05d9\u05dd \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d0\u05ea \u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05db\u05d9\u05d5\u05d5\u05df \u05d4\u05e0\u05d2\u05d3\u05d9 \u05e9\u05dc \u05d4\u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05e9\u05dc \u05de\u05d9\u05e0\u05d9-\u05d1\u05d0\u05e5'.\",\n  \"Replace\"\n);\n\n// 4) Body paragraph 2\nitems[3].insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05e0\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0 \u05e9\u05dc \u05d8\u05e7\u05e1\u05d8 \u05d4\u05de\u05db\u05d9\u05dc \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05d5\u05de\u05db\u05e0\u05d9\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05d9\u05e0\u05d9-\u05d1\u05d0\u05e5' \u05db\u05dc \u05db\u05de\u05d4 \u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc GD. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d0\u05d5 \u05db\u05de\u05d4 \u05d3\u05d1\u05e8\u05d9\u05dd \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d9\u05dd. \u05dc\u05de\u05e9\u05dc \u05db\u05de\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05d5\u05de\u05df \u05e2\u05dc\u05d9\u05d5 \u05dc\u05e4\u05e0\u05d9 \u05d4\u05ea\u05d7\u05dc\u05ea \u05d4\u05d6\u05e8\u05e7\u05ea \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05dc\u05d0 \u05de\u05e9\u05e4\u05d9\u05e2 \u05e2\u05dc \u05de\u05e1\u05e4\u05e8 \u05d4\u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d5\u05ea \u05d4\u05e0\u05d3\u05e8\u05e9 \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05e9\u05dc \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9. \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05d5\u05ea\u05e8 \u05f4\u05d9\u05d3\u05e2\u05f4 \u05d4\u05e0\u05de\u05e6\u05d0 \u05db\u05d1\u05e8 \u05d1\u05de\u05d5\u05d3\u05dc \u05dc\u05ea\u05d5\u05e8\u05dd \u05dc\u05de\u05d4\u05d9\u05e8\u05d5\u05ea \u05d4\u05dc\u05de\u05d9\u05d3\u05d4.\",\n  \"Replace\"\n);\n\n// 5) Body paragraph 3\nitems[4].insertText(\n  \"\u05e9\u05e0\u05d9\u05ea, \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e8\u05d0\u05d4 \u05e9\u05de\u05d4\u05d9\u05e8\u05d5\u05ea \u05d4\u05dc\u05de\u05d9\u05d3\u05d4 \u05e9\u05dc \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05dc\u05d0 \u05de\u05d5\u05e9\u05e4\u05e2\u05ea \u05de\u05de\u05ea\u05d9 \u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05dc\u05d4\u05d6\u05e8\u05d9\u05e7 \u05dc\u05de\u05d5\u05d3\u05dc \u05d0\u05ea \u05d4\u05d9\u05d3\u05e2. \u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05d5\u05d3\u05dc \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d0\u05d5 \u05d3\u05d5\u05d5\u05e7\u05d0 \u05ea\u05dc\u05de\u05d9\u05d3 \u05d9\u05d5\u05ea\u05e8 \u05d8\u05d5\u05d1. \u05d5\u05d9\u05e9 \u05e2\u05d5\u05d3 \u05db\u05de\u05d4 \u05ea\u05d2\u05dc\u05d9\u05d5\u05ea \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8.\",\n  \"Replace\"\n);\n\n// 6) Former link paragraph becomes a new body paragraph (4th body paragraph)\nitems[5].insertText(\n  \"\u05d0\u05d9\u05da \u05d1\u05d5\u05d3\u05e7\u05d9\u05dd \u05d4\u05d0\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05db\u05df \u05dc\u05de\u05d3 \u05d0\u05ea \u05d4\u05d9\u05d3\u05e2 \u05d4\u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05e9\u05d4\u05d6\u05e8\u05e7\u05e0\u05d5 - \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05dc\u05d0 \u05de\u05e8\u05d7\u05d9\u05d1\u05d9\u05dd \u05e2\u05dc \u05db\u05da \u05d0\u05d1\u05dc \u05db\u05e0\u05e8\u05d0\u05d4 \u05d6\u05d4 \u05de\u05d7\u05d5\u05e9\u05d1 \u05d3\u05e8\u05da likelihood \u05e9\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e0\u05db\u05d5\u05e0\u05d4 \u05e2\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05dc\u05d2\u05d1\u05d9 \u05e4\u05d9\u05e1\u05ea \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05d6\u05d4, \u05dc\u05de\u05e9\u05dc \u05f4\u05de\u05d4 \u05e2\u05d9\u05e8 \u05d4\u05d1\u05d9\u05e8\u05d4 \u05e9\u05dc \u05e6\u05e8\u05e4\u05ea\u05f4.\",\n  \"Replace\"\n);\n\nawait context.sync();\n\n// 7) New paragraph added at the end holding the new link.\nconst last = body.paragraphs.getLast();\nlast.insertParagraph(\"https://arxiv.org/abs/2406.11813\", \"After\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Date line\n$d.Paragraphs.Item(1).Range.Text = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 04.07.24:\u26a1\ufe0f\ud83d\ude80\"\n\n# 2) Title\n$d.Paragraphs.Item(2).Range.Text = \"How Do Large Language Models Acquire Factual Knowledge During Pretraining?\"\n\n# 3) Body paragraph 1\n$d.Paragraphs.Item(3).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 \u05e0\u05d5\u05e9\u05d0 \u05de\u05ea\u05d9 \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d0\u05e9\u05db\u05e8\u05d4 \u05e8\u05d5\u05db\u05e9\u05d9\u05dd \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 (\u05dc\u05de\u05e9\u05dc \u05e9\u05e2\u05d9\u05e8 \u05d1\u05d9\u05e8\u05d4 \u05e9\u05dc \u05e6\u05e8\u05e4\u05ea \u05d4\u05d9\u05d0 \u05e4\u05e8\u05d9\u05e1) \u05d1\u05de\u05d4\u05dc\u05da \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05e7\u05d3\u05d9\u05dd. \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05d2\u05dd \u05d1\u05d5\u05d3\u05e7 \u05db\u05de\u05d4 \u05d6\u05de\u05df \u05dc\u05d5\u05e7\u05d7 \u05dc\u05e9\u05db\u05d5\u05d7 \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9. \u05d0\u05d5\u05e7\u05d9\u05d9, \u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05e9\u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05e9\u05dc\u05e0\u05d5 \u05e2\u05dd \u05d0\u05d7\u05ea \u05d4\u05e6\u05d5\u05e8\u05d5\u05ea \u05e9\u05dc \u05de\u05e9\u05e4\u05d7\u05ea \u05de\u05d5\u05e8\u05d3 \u05d4\u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 (gradient descent \u05d0\u05d5 GD). \u05d1\u05d3\u05f4\u05db \u05d3\u05d5\u05d2\u05de\u05d9\u05dd \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d4\u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc\u05e0\u05d5 (\u05de\u05d9\u05e0\u05d9-\u05d1\u05d0\u05e5') \u05d5\u05de\u05d6\u05d9\u05d6\u05d9\u05dd \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d0\u05ea \u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05db\u05d9\u05d5\u05d5\u05df \u05d4\u05e0\u05d2\u05d3\u05d9 \u05e9\u05dc \u05d4\u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05e9\u05dc \u05de\u05d9\u05e0\u05d9-\u05d1\u05d0\u05e5'.\"\n\n# 4) Body paragraph 2\n$d.Paragraphs.Item(4).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05e0\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0 \u05e9\u05dc \u05d8\u05e7\u05e1\u05d8 \u05d4\u05de\u05db\u05d9\u05dc \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05d5\u05de\u05db\u05e0\u05d9\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05d9\u05e0\u05d9-\u05d1\u05d0\u05e5' \u05db\u05dc \u05db\u05de\u05d4 \u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc GD. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d0\u05d5 \u05db\u05de\u05d4 \u05d3\u05d1\u05e8\u05d9\u05dd \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d9\u05dd. \u05dc\u05de\u05e9\u05dc \u05db\u05de\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05d5\u05de\u05df \u05e2\u05dc\u05d9\u05d5 \u05dc\u05e4\u05e0\u05d9 \u05d4\u05ea\u05d7\u05dc\u05ea \u05d4\u05d6\u05e8\u05e7\u05ea \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05dc\u05d0 \u05de\u05e9\u05e4\u05d9\u05e2 \u05e2\u05dc \u05de\u05e1\u05e4\u05e8 \u05d4\u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d5\u05ea \u05d4\u05e0\u05d3\u05e8\u05e9 \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05e9\u05dc \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9. \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05d5\u05ea\u05e8 \u05f4\u05d9\u05d3\u05e2\u05f4 \u05d4\u05e0\u05de\u05e6\u05d0 \u05db\u05d1\u05e8 \u05d1\u05de\u05d5\u05d3\u05dc \u05dc\u05ea\u05d5\u05e8\u05dd \u05dc\u05de\u05d4\u05d9\u05e8\u05d5\u05ea \u05d4\u05dc\u05de\u05d9\u05d3\u05d4.\"\n\n# 5) Body paragraph 3\n$d.Paragraphs.Item(5).Range.Text = \"\u05e9\u05e0\u05d9\u05ea, \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e8\u05d0\u05d4 \u05e9\u05de\u05d4\u05d9\u05e8\u05d5\u05ea \u05d4\u05dc\u05de\u05d9\u05d3\u05d4 \u05e9\u05dc \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05dc\u05d0 \u05de\u05d5\u05e9\u05e4\u05e2\u05ea \u05de\u05de\u05ea\u05d9 \u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05dc\u05d4\u05d6\u05e8\u05d9\u05e7 \u05dc\u05de\u05d5\u05d3\u05dc \u05d0\u05ea \u05d4\u05d9\u05d3\u05e2. \u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05d5\u05d3\u05dc \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d0\u05d5 \u05d3\u05d5\u05d5\u05e7\u05d0 \u05ea\u05dc\u05de\u05d9\u05d3 \u05d9\u05d5\u05ea\u05e8 \u05d8\u05d5\u05d1. \u05d5\u05d9\u05e9 \u05e2\u05d5\u05d3 \u05db\u05de\u05d4 \u05ea\u05d2\u05dc\u05d9\u05d5\u05ea \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8.\"\n\n# 6) Former link paragraph becomes a new body paragraph (4th body paragraph)\n$d.Paragraphs.Item(6).Range.Text = \"\u05d0\u05d9\u05da \u05d1\u05d5\u05d3\u05e7\u05d9\u05dd \u05d4\u05d0\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05db\u05df \u05dc\u05de\u05d3 \u05d0\u05ea \u05d4\u05d9\u05d3\u05e2 \u05d4\u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05e9\u05d4\u05d6\u05e8\u05e7\u05e0\u05d5 - \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05dc\u05d0 \u05de\u05e8\u05d7\u05d9\u05d1\u05d9\u05dd \u05e2\u05dc \u05db\u05da \u05d0\u05d1\u05dc \u05db\u05e0\u05e8\u05d0\u05d4 \u05d6\u05d4 \u05de\u05d7\u05d5\u05e9\u05d1 \u05d3\u05e8\u05da likelihood \u05e9\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e0\u05db\u05d5\u05e0\u05d4 \u05e2\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05dc\u05d2\u05d1\u05d9 \u05e4\u05d9\u05e1\u05ea \u05d9\u05d3\u05e2 \u05e2\u05d5\u05d1\u05d3\u05ea\u05d9 \u05d6\u05d4, \u05dc\u05de\u05e9\u05dc \u05f4\u05de\u05d4 \u05e2\u05d9\u05e8 \u05d4\u05d1\u05d9\u05e8\u05d4 \u05e9\u05dc \u05e6\u05e8\u05e4\u05ea\u05f4.\"\n\n# 7) New paragraph added at the end holding the new link.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"https://arxiv.org/abs/2406.11813\"\n"}
